$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.536.89'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.851.31'
$ws.Range('E3').Value = '  +3.51%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '430.44'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.23'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.845.23'
$ws.Range('E7').Value = '  +3.57%  '
$ws.Range('E8').Value = '  -5.32%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('E10').Value = '  -5.70%  '
$ws.Range('E11').Value = '  -8.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000362'
$ws.Range('E12').Value = '  -10.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.80'
$ws.Range('E13').Value = '  -4.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.454.99'
$ws.Range('E14').Value = '  +3.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '10.03'
$ws.Range('E15').Value = '  -4.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.77'
$ws.Range('E16').Value = '  +19.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.844.50'
$ws.Range('E17').Value = '  +3.80%  '
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.56'
$ws.Range('E19').Value = '  -5.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '66.874.88'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('E21').Value = '  -6.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '407.71'
$ws.Range('E22').Value = '  -8.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.41'
$ws.Range('E23').Value = '  -12.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.09'
$ws.Range('E24').Value = '  -5.22%  '
$ws.Range('E25').Value = '  -4.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.18'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('E27').Value = '  +12.62%  '
$ws.Range('E28').Value = '  -2.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.51'
$ws.Range('E29').Value = '  -6.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '688.83'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.45'
$ws.Range('E31').Value = '  -2.81%  '
$ws.Range('E32').Value = '  -2.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.68'
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.14'
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('E35').Value = '  -8.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.75'
$ws.Range('E36').Value = '  -7.40%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0806'
$ws.Range('E37').Value = '  +7.93%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '55.27'
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.07'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  -8.00%  '
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('E44').Value = '  -8.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '148.06'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.27'
$ws.Range('E46').Value = '  -5.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.07'
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.10'
$ws.Range('E48').Value = '  -4.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '26.11'
$ws.Range('E49').Value = '  -8.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.78'
$ws.Range('E50').Value = '  -3.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.54'
$ws.Range('E51').Value = '  -4.91%  '
